# Updates the cryptos price/volume table with freshly scraped values.
# Price column D holds numeric-looking strings (e.g. "0.999", "598.42")
# that must stay TEXT (as in the source workbook) rather than be
# auto-converted to numbers by Excel's type inference. For those cells we
# assign the value with a leading apostrophe (Excel's standard "force
# text" convention) and then reset the cell style to "Normal" so no
# stray number-format/quote-prefix style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '72.447.39'
$ws.Cells.Item(2, 5).Value = '  +0.38%  '

$ws.Cells.Item(3, 4).Value = '2.664.85'
$ws.Cells.Item(3, 5).Value = '  +1.83%  '

$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

$ws.Cells.Item(5, 4).Value = "'598.42"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.83%  '

$ws.Cells.Item(6, 4).Value = "'175.49"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.43%  '

$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 4).Value = "'0.525"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Cells.Item(9, 4).Value = '2.663.44'
$ws.Cells.Item(9, 5).Value = '  +1.84%  '

$ws.Cells.Item(10, 4).Value = "'0.169"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.42%  '

$ws.Cells.Item(11, 5).Value = '  +2.35%  '

$ws.Cells.Item(12, 5).Value = '  +1.49%  '

$ws.Cells.Item(13, 5).Value = '  -0.32%  '

$ws.Cells.Item(14, 4).Value = '3.149.97'
$ws.Cells.Item(14, 5).Value = '  +1.39%  '

$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '72.233.04'
$ws.Cells.Item(15, 5).Value = '  +0.23%  '

$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(16, 4).Value = "'0.0000185"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.94%  '

$ws.Cells.Item(17, 4).Value = "'26.29"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.71%  '

$ws.Cells.Item(18, 4).Value = '2.661.88'
$ws.Cells.Item(18, 5).Value = '  +1.65%  '

$ws.Cells.Item(19, 4).Value = "'12.26"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +6.07%  '

$ws.Cells.Item(20, 4).Value = "'8.25"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +4.45%  '

$ws.Cells.Item(21, 4).Value = "'370.95"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.88%  '

$ws.Cells.Item(22, 5).Value = '  +0.39%  '

$ws.Cells.Item(23, 5).Value = '  +1.85%  '

$ws.Cells.Item(24, 4).Value = "'72.10"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.20%  '

$ws.Cells.Item(25, 5).Value = '  +0.14%  '

$ws.Cells.Item(26, 5).Value = '  -0.92%  '

$ws.Cells.Item(27, 5).Value = '  -1.25%  '

$ws.Cells.Item(28, 4).Value = '2.799.74'
$ws.Cells.Item(28, 5).Value = '  +1.77%  '

$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.04%  '

$ws.Cells.Item(30, 4).Value = '0.0₃0970'
$ws.Cells.Item(30, 5).Value = '  +2.35%  '

$ws.Cells.Item(31, 4).Value = "'8.09"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.86%  '

$ws.Cells.Item(32, 4).Value = "'501.01"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -3.44%  '

$ws.Cells.Item(33, 5).Value = '  -1.88%  '

$ws.Cells.Item(34, 5).Value = '  +0.06%  '

$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.09%  '

$ws.Cells.Item(36, 4).Value = "'162.95"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.41%  '

$ws.Cells.Item(37, 4).Value = "'19.54"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.52%  '

$ws.Cells.Item(38, 5).Value = '  +0.56%  '

$ws.Cells.Item(39, 5).Value = '  -0.58%  '

$ws.Cells.Item(40, 5).Value = '  -1.24%  '

$ws.Cells.Item(41, 4).Value = "'1.78"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -2.46%  '

$ws.Cells.Item(42, 5).Value = '  +0.05%  '

$ws.Cells.Item(43, 4).Value = "'4.99"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.81%  '

$ws.Cells.Item(44, 4).Value = "'0.333"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.63%  '

$ws.Cells.Item(45, 5).Value = '  -0.98%  '

$ws.Cells.Item(46, 4).Value = "'156.56"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +4.30%  '

$ws.Cells.Item(47, 4).Value = "'39.52"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +0.12%  '

$ws.Cells.Item(48, 5).Value = '  +1.86%  '

$ws.Cells.Item(49, 4).Value = "'0.559"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.33%  '

$ws.Cells.Item(50, 5).Value = '  +2.23%  '

$ws.Cells.Item(51, 4).Value = "'0.0755"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.31%  '
